$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 149 to hold a new week's data (2022-02-18, serial 44610).
# This pushes the existing rows 149-156 down to 152-159.
$ws.Rows("149:151").Insert()

# New row 149: Especial
$ws.Cells.Item(149, 1).Value = 12
$ws.Cells.Item(149, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(149, 3).Value = "Metropolitana"
$ws.Cells.Item(149, 4).Value = 44610
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 13
$ws.Cells.Item(149, 6).Value = 100112043
$ws.Cells.Item(149, 7).Value = "Pepino dulce"
$ws.Cells.Item(149, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(149, 9).Value = "Especial"
$ws.Cells.Item(149, 10).Value = 680
$ws.Cells.Item(149, 11).Value = 14000
$ws.Cells.Item(149, 12).Value = 15000
$ws.Cells.Item(149, 13).Value = 14588
$ws.Cells.Item(149, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(149, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(149, 16).Value = 810
$ws.Cells.Item(149, 17).Value = 18
$ws.Cells.Item(149, 18).Value = "Hortaliza"

# New row 150: Primera
$ws.Cells.Item(150, 1).Value = 12
$ws.Cells.Item(150, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44610
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = 100112043
$ws.Cells.Item(150, 7).Value = "Pepino dulce"
$ws.Cells.Item(150, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 1400
$ws.Cells.Item(150, 11).Value = 11000
$ws.Cells.Item(150, 12).Value = 12000
$ws.Cells.Item(150, 13).Value = 11514
$ws.Cells.Item(150, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(150, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(150, 16).Value = 640
$ws.Cells.Item(150, 17).Value = 18
$ws.Cells.Item(150, 18).Value = "Hortaliza"

# New row 151: Segunda
$ws.Cells.Item(151, 1).Value = 12
$ws.Cells.Item(151, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(151, 3).Value = "Metropolitana"
$ws.Cells.Item(151, 4).Value = 44610
$ws.Cells.Item(151, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 6).Value = 100112043
$ws.Cells.Item(151, 7).Value = "Pepino dulce"
$ws.Cells.Item(151, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(151, 9).Value = "Segunda"
$ws.Cells.Item(151, 10).Value = 760
$ws.Cells.Item(151, 11).Value = 8500
$ws.Cells.Item(151, 12).Value = 9000
$ws.Cells.Item(151, 13).Value = 8737
$ws.Cells.Item(151, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(151, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(151, 16).Value = 485
$ws.Cells.Item(151, 17).Value = 18
$ws.Cells.Item(151, 18).Value = "Hortaliza"
